$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()
$ws.Activate()
$ws.Select()
